# Rename worksheets to drop underscores / shorten certain labels,
# matching the updated "dalila" script naming convention.

$wb = $excel.ActiveWorkbook

$renames = @{
    "MeV.Fib_Unknown.8"        = "MeV.FibUnknown.8"
    "MeV.Endothelial_Injury.4" = "MeV.EndothelialInjury.4"
    "MeV.Immune_doublets.0"    = "MeV.ImmuneDoublets.0"
    "MeV.Low_Quality.0"        = "MeV.LowQuality.0"
    "MeV.Proliferative_Fibr.0" = "MeV.FibProlif.0"
    "MeV.Epithelial_ECad.0"    = "MeV.EpithelialECad.0"
    "MeV.Fib_CD34.7"           = "MeV.FibCD34.7"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renames.ContainsKey($oldName)) {
        $ws.Name = $renames[$oldName]
    }
}
